# Rename the "Latest Payment Date" column header on the Expense sheet to "Payment Date",
# and leave the selection/active sheet on the Expense sheet at cell E1.

$wb = $excel.ActiveWorkbook

$wsExpense = $wb.Worksheets.Item("Expense")
$wsExpense.Range("E1").Value = "Payment Date"

$wsExpense.Activate()
$wsExpense.Range("E1").Select()
